$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that used to sit right
#    after the title heading.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Berry Burst for Free - Exciting
#    Cluster Pays Slot Machine") right before the final "Prompt: ..."
#    paragraph, matching the structure of the removed meta-description
#    bold run (leading empty run + bold run).
$count = $d.Paragraphs.Count
$promptPara = $d.Paragraphs($count)
$insertPos = $promptPara.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Berry Burst for Free - Exciting Cluster Pays Slot Machine</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$insertRange.InsertXML($newParaXml)

# The insertion leaves a stray empty paragraph between the new bold
# paragraph and the (still untouched) "Prompt: ..." paragraph - remove it.
$strayParaIndex = $count + 1
$d.Paragraphs($strayParaIndex).Range.Delete()

# 3. Replace the old "Prompt: ..." text with the new meta-description
#    sentence, keeping its italic run formatting intact.
$d.Content.Find.Execute(
    "Prompt: Create a feature image for Berryburst that captures the essence of the game. The image should be in cartoon-style and feature a happy Maya warrior with glasses. The warrior should be surrounded by various fruits exploding from behind as a symbol of the game's sweet and flavorful theme. The image should convey the excitement and joy of playing Berryburst while also highlighting its fruit-inspired roots. Make sure to use bright, bold colors that pop and draw the viewer's attention. The image size should be suitable for use on online platforms such as social media, blogs, and casino websites.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover our review of Berry Burst, an innovative slot machine with cluster pays, free spins, and excellent graphics. Play now for free.",
    2)
